$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# E3: 28/5/2019 11:25 -> 27/5/2019 11:25
$ws.Range("E3").Value = "27/5/2019 11:25"

# C9: hermanliran@gmail.com -> halachme@gmail.com (value + hyperlink)
$ws.Range("C9").Value = "halachme@gmail.com"
$ws.Range("C9").Hyperlinks.Item(1).Address = "mailto:halachme@gmail.com"
$ws.Range("C9").Hyperlinks.Item(1).TextToDisplay = "halachme@gmail.com"

# Update selection to E4
$ws.Range("E4").Select()
